# Completed S3 error audit: append newly documented S3 error codes to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 43; Code = 412341; Desc = "InvalidBucketState - The request is not valid with the current state of the bucket." },
    @{ Row = 44; Code = 412342; Desc = "InvalidEncryptionAlgorithmError - The encryption request you specified is not valid. The valid value is AES256." },
    @{ Row = 45; Code = 412343; Desc = "InvalidObjectState - The operation is not valid for the current state of the object." },
    @{ Row = 46; Code = 412344; Desc = "InvalidRequest - Some part of the specified request is invalid." },
    @{ Row = 47; Code = 412345; Desc = "NoSuchLifecycleConfiguration - The lifecycle configuration does not exist." },
    @{ Row = 48; Code = 412346; Desc = "RestoreAlreadyInProgress - Object restore is already in progress." },
    @{ Row = 49; Code = 412347; Desc = "ServiceUnavailable - Reduce your request rate." }
)

foreach ($item in $newRows) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Code
    $ws.Cells.Item($item.Row, 2).Value = $item.Desc
}

# Update the view to match the final state: scrolled down and the new
# last cell (B52, three rows below the last populated data row) selected.
$ws.Range("B52").Select() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 23
$aw.ScrollColumn = 1
